# Daily crypto-symbol refresh: updates Price / Volume(1h) / Data / Hora columns
# for the coinranking.com-sourced rows, matching the GitHub Actions commit
# "Updated symbol list on Sat Dec 24 00:23:12 UTC 2022 with GitHub Actions".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay plain text even when it looks like a
# number (e.g. "246.03", "0") or it would silently become a numeric cell.
# Setting NumberFormat to "@" (Text) first, exactly like pre-formatting a cell
# as Text in the UI before typing into it, keeps the inlineStr/text storage.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Row 2
Set-TextValue "D2" '246.03'
$ws.Range("F2").Value = '24-12-2022'
Set-TextValue "G2" '0'

# Row 3
Set-TextValue "D3" '22.16'
$ws.Range("F3").Value = '24-12-2022'
Set-TextValue "G3" '0'

# Row 4
Set-TextValue "D4" '5.324'
$ws.Range("F4").Value = '24-12-2022'
Set-TextValue "G4" '0'

# Row 5
Set-TextValue "D5" '0.05882'
$ws.Range("F5").Value = '24-12-2022'
Set-TextValue "G5" '0'

# Row 6
Set-TextValue "D6" '3.397'
$ws.Range("F6").Value = '24-12-2022'
Set-TextValue "G6" '0'

# Row 7
Set-TextValue "D7" '6.379'
$ws.Range("F7").Value = '24-12-2022'
Set-TextValue "G7" '0'

# Row 8
Set-TextValue "D8" '0.8117'
$ws.Range("F8").Value = '24-12-2022'
Set-TextValue "G8" '0'

# Row 9
Set-TextValue "D9" '0.9586'
$ws.Range("F9").Value = '24-12-2022'
Set-TextValue "G9" '0'

# Row 10
Set-TextValue "D10" '0.1416'
$ws.Range("F10").Value = '24-12-2022'
Set-TextValue "G10" '0'

# Row 11
Set-TextValue "D11" '0.03499'
$ws.Range("E11").Value = '10LiechtensteinCryptoassetsExchangeLCXBestin24h'
$ws.Range("F11").Value = '24-12-2022'
Set-TextValue "G11" '0'

# Row 12
Set-TextValue "D12" '0.07374'
$ws.Range("F12").Value = '24-12-2022'
Set-TextValue "G12" '0'

# Row 13
Set-TextValue "D13" '0.03043'
$ws.Range("F13").Value = '24-12-2022'
Set-TextValue "G13" '0'

# Row 14
Set-TextValue "D14" '4.438'
$ws.Range("F14").Value = '24-12-2022'
Set-TextValue "G14" '0'

# Row 15
Set-TextValue "D15" '0.09389'
$ws.Range("F15").Value = '24-12-2022'
Set-TextValue "G15" '0'

# Row 16
Set-TextValue "D16" '0.001598'
$ws.Range("F16").Value = '24-12-2022'
Set-TextValue "G16" '0'

# Row 17
Set-TextValue "D17" '0.04826'
$ws.Range("F17").Value = '24-12-2022'
Set-TextValue "G17" '0'

# Row 18
Set-TextValue "D18" '0.0005901'
$ws.Range("F18").Value = '24-12-2022'
Set-TextValue "G18" '0'

# Row 19
Set-TextValue "D19" '0.006024'
$ws.Range("F19").Value = '24-12-2022'
Set-TextValue "G19" '0'

# Row 20
Set-TextValue "D20" '0.004088'
$ws.Range("F20").Value = '24-12-2022'
Set-TextValue "G20" '0'

# Row 21
Set-TextValue "D21" '0.0009852'
$ws.Range("F21").Value = '24-12-2022'
Set-TextValue "G21" '0'

# Row 22
Set-TextValue "D22" '0.00009703'
$ws.Range("F22").Value = '24-12-2022'
Set-TextValue "G22" '0'

# Row 23
Set-TextValue "D23" '3.691'
$ws.Range("F23").Value = '24-12-2022'
Set-TextValue "G23" '0'

# Row 24
$ws.Range("F24").Value = '24-12-2022'
Set-TextValue "G24" '0'

# Row 25
$ws.Range("F25").Value = '24-12-2022'
Set-TextValue "G25" '0'

# Row 26
Set-TextValue "D26" '0.1310'
$ws.Range("F26").Value = '24-12-2022'
Set-TextValue "G26" '0'

# Row 27
Set-TextValue "D27" '0.0002472'
$ws.Range("F27").Value = '24-12-2022'
Set-TextValue "G27" '0'

# Row 28
$ws.Range("F28").Value = '24-12-2022'
Set-TextValue "G28" '0'

# Row 29
$ws.Range("F29").Value = '24-12-2022'
Set-TextValue "G29" '0'

# Row 30
$ws.Range("F30").Value = '24-12-2022'
Set-TextValue "G30" '0'

# Row 31
$ws.Range("F31").Value = '24-12-2022'
Set-TextValue "G31" '0'

# Row 32
$ws.Range("F32").Value = '24-12-2022'
Set-TextValue "G32" '0'

# Row 33
$ws.Range("F33").Value = '24-12-2022'
Set-TextValue "G33" '0'

# Row 34
$ws.Range("F34").Value = '24-12-2022'
Set-TextValue "G34" '0'

# Row 35
$ws.Range("F35").Value = '24-12-2022'
Set-TextValue "G35" '0'

# Row 36
$ws.Range("F36").Value = '24-12-2022'
Set-TextValue "G36" '0'

# Row 37
$ws.Range("F37").Value = '24-12-2022'
Set-TextValue "G37" '0'

# Row 38
$ws.Range("F38").Value = '24-12-2022'
Set-TextValue "G38" '0'

# Row 39
$ws.Range("F39").Value = '24-12-2022'
Set-TextValue "G39" '0'

# Row 40
Set-TextValue "D40" '0.03906'
$ws.Range("F40").Value = '24-12-2022'
Set-TextValue "G40" '0'

# Row 41
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue "D41" '0.1073'
$ws.Range("E41").Value = '40BKEXTokenBKK'
$ws.Range("F41").Value = '24-12-2022'
Set-TextValue "G41" '0'

# Row 42
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextValue "D42" '0.003001'
$ws.Range("E42").Value = '41CEJICEJI'
$ws.Range("F42").Value = '24-12-2022'
Set-TextValue "G42" '0'

# Row 43
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue "D43" '0.006742'
$ws.Range("E43").Value = '42KickTokenKICK'
$ws.Range("F43").Value = '24-12-2022'
Set-TextValue "G43" '0'

# Row 44
Set-TextValue "D44" '0.005896'
$ws.Range("F44").Value = '24-12-2022'
Set-TextValue "G44" '0'

# Row 45
Set-TextValue "D45" '0.00005691'
$ws.Range("F45").Value = '24-12-2022'
Set-TextValue "G45" '0'

# Row 46
Set-TextValue "D46" '0.00000000750'
$ws.Range("F46").Value = '24-12-2022'
Set-TextValue "G46" '0'

# Row 47
Set-TextValue "D47" '0.6002'
$ws.Range("E47").Value = '46CoinbaseStockTokenCOIN'
$ws.Range("F47").Value = '24-12-2022'
Set-TextValue "G47" '0'

# Row 48
Set-TextValue "D48" '0.05212'
$ws.Range("F48").Value = '24-12-2022'
Set-TextValue "G48" '0'

# Row 49
$ws.Range("F49").Value = '24-12-2022'
Set-TextValue "G49" '0'

# Row 50
Set-TextValue "D50" '0.01010'
$ws.Range("F50").Value = '24-12-2022'
Set-TextValue "G50" '0'

# Row 51
$ws.Range("F51").Value = '24-12-2022'
Set-TextValue "G51" '0'
